# The "RegisterUserData" sheet previously held a wide table (Name, Email,
# Password, Interest, Gender, State, Hobby) with two sample rows. Per the
# commit, this is trimmed down to just two columns (Name, State) with a
# single data row ("Darryal" / "Goa"), leaving the sheet as a small
# 2x2 (A1:B2) block for TDD purposes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegisterUserData")

# Remove the hyperlinks that live on column C (Password column); that
# column is going away entirely.
$ws.Hyperlinks.Delete()

# Clear the custom column widths/styles on columns A and B so no <cols>
# override remains (the target sheet just uses sheet defaults).
$ws.Columns.Item(1).ClearFormats()
$ws.Columns.Item(2).ClearFormats()

# Drop columns C:G (Password, Interest, Gender, State, Hobby) - data shifts
# left, so the old "State" column (F) ends up gone too; we re-add a new
# State column below as column B.
$ws.Range("C1:G3").EntireColumn.Delete()

# Drop the third data row, leaving only the header row and one data row.
$ws.Range("A3:B3").EntireRow.Delete()

# Write the new, smaller table.
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "State"
$ws.Range("A2").Value = "Darryal"
$ws.Range("B2").Value = "Goa"

# Let the remaining two rows fall back to the sheet's default row height
# instead of keeping the old explicit height.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(2).AutoFit()

# Match the author's cursor position recorded in the saved file.
[void]$ws.Range("C2").Select()
